# Integrate frontend and implement 'UPDATE_OWNER_SCREEN'
#
# The catalog sheet is reshuffled: the old Admin-screen block (rows 5-7,
# using SHOW_ACTIVE_ADMIN_SCREEN / SHOW_INACTIVE_ADMIN_SCREEN) is replaced
# by a single SHOW_OWNER_SCREEN row, and the Login block's
# SHOW_INVALID_INPUT_ERROR_MESSAGE / SHOW_LOGIN_SCREEN rows shift up to
# fill the gap. A new blank, bordered block (rows 12-15) is appended below
# the existing Owner/UPDATE_OWNER_SCREEEN row, and that Owner row (11) gets
# a distinct "open box" border (no bottom edge) to visually group it with
# the new blank rows beneath it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content updates ------------------------------------------------

# Former "Login / Admin" rows 5-7 collapse: row 5 becomes the new
# SHOW_OWNER_SCREEN entry (no Parameters/Comment), and the
# SHOW_INVALID_INPUT_ERROR_MESSAGE / SHOW_LOGIN_SCREEN pair shifts from
# rows 6-7 up into rows 6-7 in swapped order (Invalid Input now first).
$ws.Range("B5").Value = "SHOW_OWNER_SCREEN"
$ws.Range("C5").ClearContents()
$ws.Range("D5").ClearContents()

$ws.Range("B6").Value = "SHOW_INVALID_INPUT_ERROR_MESSAGE"
$ws.Range("C6").Value = "msg"

$ws.Range("B7").Value = "SHOW_LOGIN_SCREEN"
$ws.Range("C7").ClearContents()

# --- New blank rows 12-15, bordered like the rest of the table -----------

$srcBlank = $ws.Range("A8:D8")
$dstBlank = $ws.Range("A12:D15")
$srcBlank.Copy()
$dstBlank.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 11 (Owner / UPDATE_OWNER_SCREEEN) gets an "open box" border -----
# (left/right/top thin, no bottom) instead of the full box used elsewhere.

$srcFull = $ws.Range("A2:D2")
$dstRow11 = $ws.Range("A11:D11")
$srcFull.Copy()
$dstRow11.PasteSpecial(-4122)
$excel.CutCopyMode = $false
$dstRow11.Borders.Item(9).LineStyle = -4142
$ws.Range("A11").Font.Bold = $true

# --- View settings ---------------------------------------------------------

$excel.ActiveWindow.Zoom = 115
$ws.Range("C11").Select()

Write-Host "edit complete"
